# Update menu items: new dishes added (Schezwan Sandwich, Club Sandwich,
# Steamed Veg Momos, Kurkure Paneer Momos), Pizza rows moved and prices
# updated, and the now-obsolete "Half" price column cleared for rows
# that no longer offer a half portion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Chole Bhature"
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = "Chole Bhature.jpeg"

$ws.Range("A3").Value = "Chole chawal"
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = "Chole Chawal.jpeg"

$ws.Range("A4").Value = "Chumin "
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = "Chumin.jpeg"

$ws.Range("A5").Value = "Chumin Paneer"
$ws.Range("B5").Value = 35
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = "Chumin Paneer.jpeg"

$ws.Range("A6").Value = "Chill Patato"
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = "chill Patato.jpeg"

$ws.Range("A7").Value = "Honey Chill patato"
$ws.Range("B7").Value = 35
$ws.Range("C7").Value = 70
$ws.Range("D7").Value = "Honey Chill patato.jpeg"

$ws.Range("A8").Value = "Aloo tikki Burger "
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = "Burger Aloo tikki.png"

$ws.Range("A9").Value = "Aloo tikki chees Burger"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = "XL Burger Aloo tikki chees.jpeg"

$ws.Range("A10").Value = "Paneer tikki chees Burger"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 70
$ws.Range("D10").Value = "Brioche Bun Burger Paneer tikki chees.jpeg"

$ws.Range("A11").Value = "Tortila Wrap Paneer"
$ws.Range("B11").Value = 50
$ws.Range("C11").Value = 60
$ws.Range("D11").Value = "Tortila Wrap Paneer.jpeg"

$ws.Range("A12").Value = "Burrito wrap Paneer"
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = 70
$ws.Range("D12").Value = "Burrito wrap Paneer.jpeg"

$ws.Range("A13").Value = "Quesadilla Paneer"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 60
$ws.Range("D13").Value = "Quesadilla Paneer.jpeg"

$ws.Range("A14").Value = "Schezwan Grilled Sandwich – Indo-Chinese fusion with Schezwan sauce, veggies, and cheese."
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = 80
$ws.Range("D14").Value = "Schezwan Grilled Sandwich.png"

$ws.Range("A15").Value = "Club Sandwich (Indian Style) – Multi-layered with veggies, green chutney, Paneer patty, and cheese"
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = 100
$ws.Range("D15").Value = "Multi-layered with veggies Paneer.png"

$ws.Range("A16").Value = "Veg Thail(Chole, Mix Veg, Rice, Raita, salad, 2 Roti, Gulab Zamun)"
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 150
$ws.Range("D16").Value = "Veg Thali(Chole, Mix Veg, Raita, salad, 4 Roti, Gulab Jamun).png"

$ws.Range("A17").Value = "Veg Special Thail(Sabzi Paneer, Mix Veg, Rice, Raita, salad, 2 Roti, Gulab Zamun)"
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 175
$ws.Range("D17").Value = "Veg Special Thail(Sabzi Paneer, Mix Veg, Raita, salad, 4 Roti, Gulab Zamun).png"

$ws.Range("A18").Value = "Veg Biryani Soya with Garlic Mayo Dip *1"
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = 100
$ws.Range("D18").Value = "Veg Biryani Soya with Garlic Mayo Dip 1.png"

$ws.Range("A19").Value = "Veg Biryani Paneer with Dip Garlic Mayo Dip *1"
$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = 150
$ws.Range("D19").Value = "Veg Biryani Paneer with Dip Garlic Mayo Dip 1.png"

$ws.Range("A20").Value = "Pizza Margireta"
$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = 80
$ws.Range("D20").Value = "Pizza Margireta.jpeg"

$ws.Range("A21").Value = "Pizza Veggi Panner( Onion and Capcium and corn)"
$ws.Range("B21").ClearContents()
$ws.Range("C21").Value = 120
$ws.Range("D21").Value = "Pizza Veggi Panner( Onion and Capcium and corn).jpeg"

$ws.Range("A22").Value = "Pizza Onion and Capcium and corn"
$ws.Range("B22").ClearContents()
$ws.Range("C22").Value = 100
$ws.Range("D22").Value = "Pizza Onion and Capcium and corn.jpeg"

$ws.Range("A23").Value = "Steamed Veg Momos/8 piece"
$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = 80
$ws.Range("D23").Value = "Steam Momos paneer.jpeg"

$ws.Range("A24").Value = "Kurkure Paneer Momo's/8 piece"
$ws.Range("B24").ClearContents()
$ws.Range("C24").Value = 120
$ws.Range("D24").Value = "Kurkure Paneer Momo's.jpeg"

$ws.Range("D28").Select()
